$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 into the new header cells I1:J1, then set their text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-15
$data = @(
    @(7, 7),
    @(5, 5),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(7, 8),
    @(7, 7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
